$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a number by Excel
# must be forced to stay text (matching the original inlineStr/text cells).
$ws.Range('D2').Value = '29.653.45'
$ws.Range('E2').Value = '  +3.55%  '
$ws.Range('D3').Value = '1.609.25'
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '212.79'
$ws.Range('E6').Value = '  +1.49%  '
$ws.Range('E7').Value = '  +0.21%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '27.22'
$ws.Range('E8').Value = '  +9.69%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '43.68'
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('E11').Value = '  +2.49%  '
$ws.Range('E12').Value = '  +1.52%  '
$ws.Range('D13').Value = '1.840.50'
$ws.Range('E13').Value = '  +2.89%  '
$ws.Range('D14').Value = '1.595.01'
$ws.Range('E14').Value = '  +1.92%  '
$ws.Range('D15').Value = '29.662.88'
$ws.Range('E15').Value = '  +3.45%  '
$ws.Range('E16').Value = '  +4.31%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '63.59'
$ws.Range('E18').Value = '  +3.38%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '240.68'
$ws.Range('E19').Value = '  +5.83%  '
$ws.Range('E20').Value = '  +3.96%  '
$ws.Range('D21').Value = '0.0₃0694'
$ws.Range('E21').Value = '  +1.97%  '
$ws.Range('E22').Value = '  +0.28%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.00'
$ws.Range('E23').Value = '  +2.12%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.23'
$ws.Range('E24').Value = '  +2.05%  '
$ws.Range('E25').Value = '  +0.95%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '155.28'
$ws.Range('E26').Value = '  +2.53%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '15.35'
$ws.Range('E27').Value = '  +4.02%  '
$ws.Range('E28').Value = '  +1.51%  '
$ws.Range('E29').Value = '  +2.96%  '
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('E31').Value = '  +3.92%  '
$ws.Range('E32').Value = '  +0.74%  '
$ws.Range('B34').Value = 'Maker'
$ws.Range('C34').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D34').Value = '1.434.22'
$ws.Range('E34').Value = '  +2.16%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.13'
$ws.Range('E35').Value = '  +4.30%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.03'
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.54'
$ws.Range('E37').Value = '  +5.05%  '
$ws.Range('E38').Value = '  +5.43%  '
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0166'
$ws.Range('E40').Value = '  +2.01%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.541'
$ws.Range('E41').Value = '  +4.75%  '
$ws.Range('E42').Value = '  +1.95%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '54.29'
$ws.Range('E43').Value = '  +27.47%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0490'
$ws.Range('E44').Value = '  +5.62%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.801'
$ws.Range('E45').Value = '  +4.59%  '
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '66.06'
$ws.Range('E47').Value = '  +3.38%  '
$ws.Range('E48').Value = '  +1.66%  '
$ws.Range('D49').Value = '1.750.42'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.927'
$ws.Range('E50').Value = '  +8.06%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '86.92'
$ws.Range('E51').Value = '  +2.55%  '
